{"js": "// The document consists of a single (otherwise empty) paragraph that only\n// contains the `_GoBack` bookmark. We need to insert a new run of text\n// right before that bookmark, with run formatting (sz/szCs/lang) matching\n// the paragraph mark's rPr.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Target the last paragraph in the body (the one being edited in the diff).\nconst targetParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst text =\n  \"Actually In the question we just have to define the appropriate Model Library, Data Access Layer Library(DAL) and Business Logic(BL) interfaces for Doctor, patient, and Appointment including service methods. So No OUTPUTS for this assignment.\";\n\n// Build a flat-OPC OOXML fragment for a single run carrying the desired\n// run properties, and insert it at the very start of the paragraph (i.e.\n// before the bookmark start/end that already live there).\nconst runOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:rPr><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:lang w:val=\"en-US\"/></w:rPr>' +\n  \"<w:t>\" + text + \"</w:t>\" +\n  \"</w:r></w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntargetParagraph.insertOoxml(runOoxml, \"Start\");\nawait context.sync();\n", "ps1": "# The document is a single (otherwise empty) paragraph that only contains\n# the `_GoBack` bookmark. Insert a new run of text right before that\n# bookmark, with run formatting (sz/szCs/lang) matching the paragraph\n# mark's rPr (sz=32, szCs=32, lang=en-US).\n$d = $word.ActiveDocument\n\n# Target the last paragraph in the document (the one being edited).\n$p = $d.Paragraphs($d.Paragraphs.Count)\n$r = $p.Range\n\n# Collapse to the very start of the paragraph so the new content lands\n# before the existing bookmarkStart/bookmarkEnd.\n$r.Collapse(1)\n\n$text = \"Actually In the question we just have to define the appropriate Model Library, Data Access Layer Library(DAL) and Business Logic(BL) interfaces for Doctor, patient, and Appointment including service methods. So No OUTPUTS for this assignment.\"\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:rPr><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:lang w:val=\"en-US\"/></w:rPr>' +\n  '<w:t>' + $text + '</w:t>' +\n  '</w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($xml)\n"}
